# EZ-3083 imail collection fixes
#
# The "@Address5@" address line is removed (merged away) and the
# paragraph that used to hold it now starts with a fresh "_GoBack"
# bookmark immediately followed by the "@Postcode@" run. Word keeps
# only a single "_GoBack" bookmark per document, so adding the new one
# automatically relocates (removes) the old one that used to sit in
# front of the "THIS NOTICE SHOULD INCLUDE..." paragraph further down.

$d = $word.ActiveDocument

# Locate the paragraph that contains the literal "@Address5@" merge
# field text.
$addr5Para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^@Address5@\r?$") {
        $addr5Para = $p
        break
    }
}

if ($addr5Para -ne $null) {
    # Delete the "@Address5@" text together with its trailing paragraph
    # mark. This merges the (now empty) paragraph into the next one
    # ("@Postcode@"), leaving a single paragraph behind.
    $mergeRange = $d.Range($addr5Para.Range.Start, $addr5Para.Range.End)
    $mergeRange.Delete()
}

# Re-find the paragraph that now holds "@Postcode@" (it is the merged
# paragraph) and drop a "_GoBack" bookmark right at its start, ahead of
# the run.
$postcodePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^@Postcode@\r?$") {
        $postcodePara = $p
        break
    }
}

if ($postcodePara -ne $null) {
    $bmStart = $postcodePara.Range.Start
    $bmRange = $d.Range($bmStart, $bmStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
